$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update category data in rows 2-4 (column A first, then column B,
# to match the shared-string insertion order)
$ws.Range("A2").Value = "KTG013"
$ws.Range("A3").Value = "KTG014"
$ws.Range("A4").Value = "KTG015"
$ws.Range("B2").Value = "Buku Komik"
$ws.Range("B3").Value = "Buku Sejarah"
$ws.Range("B4").Value = "Buku Tulis"

# Update selection to C8
$ws.Range("C8").Select()
